$d = $word.ActiveDocument

# --- Transform paragraph: "FLOAT" -> "REAL" (both occurrences: the data type
# name and the Postgres variable type name). Use Replace-All so both hits
# are updated in a single Execute call.
$d.Content.Find.Execute("FLOAT", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "REAL", 2)

# --- Closing paragraph: "2000 to 2010" -> "2000 to 2016" so it agrees with
# the "year 2016" table referenced earlier in the report.
$d.Content.Find.Execute("2000 to 2010", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2000 to 2016", 2)

# --- Relocate the automatic "_GoBack" bookmark (Word drops this at the
# location of the most recent edit) to sit right after the second "REAL",
# immediately before " variable. Below are some snippets...". Adding a
# bookmark with the reserved name "_GoBack" replaces/moves any existing
# one, removing it from its old spot at the end of the document.
$range = $d.Content
$range.Find.Execute(" variable. Below are some snippets of the code I wrote to transform the data.")
$range.Collapse(1)
$d.Bookmarks.Add("_GoBack", $range)
